$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$acctFmt = "_-* #,##0.00\ [`$€-C0A]_-;\-* #,##0.00\ [`$€-C0A]_-;_-* ""-""??\ [`$€-C0A]_-;_-@_-"

# --- New labels / values for the "Capital aportado" mini balance-sheet block (H:K) ---
# (shared-string table insertion order matches the original authoring order)
$ws.Range("H15").Value = "1 Capital aportado"
$ws.Range("I15").Value = 50000
$ws.Range("I15").NumberFormat = $acctFmt

$ws.Range("H16").Value = "Aplicaiones del capital"
$ws.Range("J16").Value = "Origigenes del Capital"

$ws.Range("J21").Value = "Capital Social"
$ws.Range("K21").Value = 50000
$ws.Range("K21").NumberFormat = $acctFmt

$ws.Range("H17").Value = "Caja"
$ws.Range("I17").Value = 50000
$ws.Range("I17").NumberFormat = $acctFmt

$ws.Range("J20").Value = "Prestamos Bancarios"
$ws.Range("K20").Value = 20000
$ws.Range("K20").NumberFormat = $acctFmt

$ws.Range("H18").Value = "+"
$ws.Range("I18").Value = 20000
$ws.Range("I18").NumberFormat = $acctFmt

# --- The rest of the I / K accounting-format column is filled but left blank ---
$blankFmtCells = @("K14","I16","K15","K16","I19","K17","K18","K19","I20","I21","I22","K22","I23","K23","I24","K24","I25","K25","I26","K26","I27","K27","I28","K28")
foreach ($addr in $blankFmtCells) {
    $ws.Range($addr).NumberFormat = $acctFmt
}

# --- Column widths for the new columns ---
$ws.Columns("H").ColumnWidth = 20.666666666666668
$ws.Columns("I").ColumnWidth = 11.166666666666666
$ws.Columns("J").ColumnWidth = 19.666666666666668
$ws.Columns("K").ColumnWidth = 11.166666666666666

# --- Picture moved/resized (side effect of the new columns) ---
$shp = $ws.Shapes.Item(1)
$shp.Left = 58.4375
$shp.Top = 15.0
$shp.Width = 494.6875
$shp.Height = 156.33732283464568

# --- Selection / scroll position left by the editing session ---
$ws.Range("H22").Select() | Out-Null
